# Slovakia 2-liga 2023-2024 — re-sync odds sheet with latest scrape.
# A brand-new match (Myjava vs Puchov, 2023-11-09 18:00) was added by the
# scraper. Because several kick-offs share an identical date/time stamp,
# the rows that tie on that timestamp get re-emitted in a new (arbitrary)
# order on every run; this script reproduces that exact re-shuffle plus
# the new row appended at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 92 (same kickoff block as rows 92-94) ----
$ws.Range("F92").Value = "Puchov"
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = "Komarno"
$ws.Range("J92").Value = 2.69
$ws.Range("L92").Value = 2.36
$ws.Range("M92").Value = "14/10/2023 14:52"
$ws.Range("N92").Value = 3.23
$ws.Range("P92").Value = 3.38
$ws.Range("Q92").Value = "14/10/2023 14:52"
$ws.Range("R92").Value = 2.31
$ws.Range("T92").Value = 2.84
$ws.Range("U92").Value = "14/10/2023 14:52"
$ws.Range("V92").Value = "https://www.betexplorer.com/football/slovakia/2-liga/msk-puchov-komarno/8dfbsaxo/"

# ---- Row 93 ----
$ws.Range("F93").Value = "Povazska Bystrica"
$ws.Range("H93").Value = "Zilina B"
$ws.Range("J93").Value = 1.77
$ws.Range("L93").Value = 2.15
$ws.Range("M93").Value = "14/10/2023 14:56"
$ws.Range("N93").Value = 3.81
$ws.Range("P93").Value = 3.74
$ws.Range("Q93").Value = "14/10/2023 14:56"
$ws.Range("R93").Value = 3.47
$ws.Range("T93").Value = 2.95
$ws.Range("U93").Value = "14/10/2023 14:56"
$ws.Range("V93").Value = "https://www.betexplorer.com/football/slovakia/2-liga/povazska-bystrica-zilina/dCmt6rFo/"

# ---- Row 94 ----
$ws.Range("F94").Value = "D. Kubin"
$ws.Range("G94").Value = 2
$ws.Range("H94").Value = "Spisska Nova Ves"
$ws.Range("J94").Value = 2
$ws.Range("L94").Value = 2.45
$ws.Range("M94").Value = "14/10/2023 14:48"
$ws.Range("N94").Value = 3.36
$ws.Range("P94").Value = 3.49
$ws.Range("Q94").Value = "14/10/2023 14:51"
$ws.Range("R94").Value = 3.15
$ws.Range("T94").Value = 2.65
$ws.Range("U94").Value = "14/10/2023 14:48"
$ws.Range("V94").Value = "https://www.betexplorer.com/football/slovakia/2-liga/d-kubin-spisska-nova-ves/IRyk4Mqb/"

# ---- Row 100 (block 100-102) ----
$ws.Range("F100").Value = "Malzenice"
$ws.Range("G100").Value = 2
$ws.Range("H100").Value = "Myjava"
$ws.Range("J100").Value = 2.59
$ws.Range("L100").Value = 3.04
$ws.Range("M100").Value = "21/10/2023 14:28"
$ws.Range("N100").Value = 3.23
$ws.Range("P100").Value = 3.14
$ws.Range("Q100").Value = "21/10/2023 14:28"
$ws.Range("R100").Value = 2.4
$ws.Range("T100").Value = 2.36
$ws.Range("U100").Value = "21/10/2023 14:23"
$ws.Range("V100").Value = "https://www.betexplorer.com/football/slovakia/2-liga/malzenice-myjava/4bSOGaT3/"

# ---- Row 101 ----
$ws.Range("F101").Value = "Presov"
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = "Petrzalka"
$ws.Range("I101").Value = 1
$ws.Range("J101").Value = 1.56
$ws.Range("L101").Value = 1.85
$ws.Range("M101").Value = "21/10/2023 14:19"
$ws.Range("N101").Value = 3.84
$ws.Range("P101").Value = 3.74
$ws.Range("Q101").Value = "21/10/2023 14:19"
$ws.Range("R101").Value = 4.69
$ws.Range("T101").Value = 3.82
$ws.Range("U101").Value = "21/10/2023 14:19"
$ws.Range("V101").Value = "https://www.betexplorer.com/football/slovakia/2-liga/presov-petrzalka/Opf2abbT/"

# ---- Row 102 ----
$ws.Range("F102").Value = "Spisska Nova Ves"
$ws.Range("H102").Value = "FK Humenne"
$ws.Range("I102").Value = 3
$ws.Range("J102").Value = 3.02
$ws.Range("L102").Value = 3.43
$ws.Range("M102").Value = "21/10/2023 14:27"
$ws.Range("N102").Value = 3.19
$ws.Range("P102").Value = 3.47
$ws.Range("Q102").Value = "21/10/2023 14:27"
$ws.Range("R102").Value = 2.13
$ws.Range("T102").Value = 2.04
$ws.Range("U102").Value = "21/10/2023 14:27"
$ws.Range("V102").Value = "https://www.betexplorer.com/football/slovakia/2-liga/spisska-nova-ves-fk-humenne/fTpYgxMj/"

# ---- Row 105 (block 105-107) ----
$ws.Range("F105").Value = "Povazska Bystrica"
$ws.Range("G105").Value = 3
$ws.Range("H105").Value = "Slovan Bratislava B"
$ws.Range("J105").Value = 1.8
$ws.Range("L105").Value = 1.46
$ws.Range("M105").Value = "28/10/2023 13:57"
$ws.Range("N105").Value = 3.56
$ws.Range("P105").Value = 4.44
$ws.Range("Q105").Value = "28/10/2023 13:57"
$ws.Range("R105").Value = 3.58
$ws.Range("T105").Value = 6.15
$ws.Range("U105").Value = "28/10/2023 13:57"
$ws.Range("V105").Value = "https://www.betexplorer.com/football/slovakia/2-liga/povazska-bystrica-slovan-bratislava/ObZzEcDM/"

# ---- Row 106 ----
$ws.Range("F106").Value = "Puchov"
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = "Pohronie"
$ws.Range("I106").Value = 4
$ws.Range("J106").Value = 1.75
$ws.Range("L106").Value = 1.71
$ws.Range("M106").Value = "28/10/2023 14:21"
$ws.Range("N106").Value = 3.6
$ws.Range("P106").Value = 3.81
$ws.Range("Q106").Value = "28/10/2023 14:21"
$ws.Range("R106").Value = 3.77
$ws.Range("T106").Value = 4.45
$ws.Range("U106").Value = "28/10/2023 14:21"
$ws.Range("V106").Value = "https://www.betexplorer.com/football/slovakia/2-liga/msk-puchov-fk-pohronie/QDGXEwbG/"

# ---- Row 107 ----
$ws.Range("F107").Value = "FK Humenne"
$ws.Range("G107").Value = 4
$ws.Range("H107").Value = "L. Mikulas"
$ws.Range("I107").Value = 1
$ws.Range("J107").Value = 2.3
$ws.Range("L107").Value = 1.96
$ws.Range("N107").Value = 3.18
$ws.Range("P107").Value = 3.59
$ws.Range("R107").Value = 2.74
$ws.Range("T107").Value = 3.52
$ws.Range("V107").Value = "https://www.betexplorer.com/football/slovakia/2-liga/fk-humenne-l-mikulas/x6FB7ejj/"

# ---- Row 109 (block 109-111) ----
$ws.Range("F109").Value = "Samorin"
$ws.Range("G109").Value = 2
$ws.Range("H109").Value = "Myjava"
$ws.Range("I109").Value = 2
$ws.Range("J109").Value = 2.32
$ws.Range("L109").Value = 2.78
$ws.Range("M109").Value = "29/10/2023 09:48"
$ws.Range("N109").Value = 3.35
$ws.Range("P109").Value = 3.53
$ws.Range("Q109").Value = "29/10/2023 09:48"
$ws.Range("R109").Value = 2.66
$ws.Range("T109").Value = 2.33
$ws.Range("U109").Value = "29/10/2023 09:48"
$ws.Range("V109").Value = "https://www.betexplorer.com/football/slovakia/2-liga/samorin-myjava/Ua2cVbrc/"

# ---- Row 110 ----
$ws.Range("F110").Value = "Petrzalka"
$ws.Range("H110").Value = "Komarno"
$ws.Range("J110").Value = 2.89
$ws.Range("K110").Value = "27/10/2023 23:42"
$ws.Range("L110").Value = 2.91
$ws.Range("M110").Value = "29/10/2023 10:02"
$ws.Range("N110").Value = 3.2
$ws.Range("O110").Value = "27/10/2023 23:42"
$ws.Range("P110").Value = 3.44
$ws.Range("Q110").Value = "29/10/2023 10:02"
$ws.Range("R110").Value = 2.2
$ws.Range("S110").Value = "27/10/2023 23:42"
$ws.Range("T110").Value = 2.29
$ws.Range("U110").Value = "29/10/2023 10:02"
$ws.Range("V110").Value = "https://www.betexplorer.com/football/slovakia/2-liga/petrzalka-komarno/GdEF6F5d/"

# ---- Row 111 ----
$ws.Range("F111").Value = "Malzenice"
$ws.Range("G111").Value = 3
$ws.Range("H111").Value = "Spisska Nova Ves"
$ws.Range("I111").Value = 1
$ws.Range("J111").Value = 1.94
$ws.Range("K111").Value = "28/10/2023 08:13"
$ws.Range("L111").Value = 2.03
$ws.Range("M111").Value = "29/10/2023 10:21"
$ws.Range("N111").Value = 3.34
$ws.Range("O111").Value = "28/10/2023 08:13"
$ws.Range("P111").Value = 3.33
$ws.Range("Q111").Value = "29/10/2023 10:29"
$ws.Range("R111").Value = 3.32
$ws.Range("S111").Value = "28/10/2023 08:13"
$ws.Range("T111").Value = 3.59
$ws.Range("U111").Value = "29/10/2023 10:21"
$ws.Range("V111").Value = "https://www.betexplorer.com/football/slovakia/2-liga/malzenice-spisska-nova-ves/KUK68yyp/"

# ---- Row 123 (block 123-125, incl. the new row) ----
$ws.Range("F123").Value = "Komarno"
$ws.Range("G123").Value = 1
$ws.Range("H123").Value = "FK Humenne"
$ws.Range("I123").Value = 1
$ws.Range("J123").Value = 1.56
$ws.Range("L123").Value = 1.67
$ws.Range("M123").Value = "11/11/2023 12:45"
$ws.Range("N123").Value = 3.77
$ws.Range("P123").Value = 3.74
$ws.Range("Q123").Value = "11/11/2023 12:45"
$ws.Range("R123").Value = 4.81
$ws.Range("T123").Value = 4.9
$ws.Range("U123").Value = "11/11/2023 12:45"
$ws.Range("V123").Value = "https://www.betexplorer.com/football/slovakia/2-liga/komarno-fk-humenne/buGqsds9/"

# ---- Row 124 ----
$ws.Range("F124").Value = "Presov"
$ws.Range("G124").Value = 2
$ws.Range("H124").Value = "Spisska Nova Ves"
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 1.18
$ws.Range("L124").Value = 1.17
$ws.Range("M124").Value = "11/11/2023 12:57"
$ws.Range("N124").Value = 6.08
$ws.Range("P124").Value = 6.84
$ws.Range("Q124").Value = "11/11/2023 12:57"
$ws.Range("R124").Value = 9.970000000000001
$ws.Range("T124").Value = 16.39
$ws.Range("U124").Value = "11/11/2023 12:57"
$ws.Range("V124").Value = "https://www.betexplorer.com/football/slovakia/2-liga/presov-spisska-nova-ves/4xJiuzCL/"

# ---- Row 125 (brand new match row, appended at the end) ----
# Clone row 124's cell formatting (bold/border style on col A, datetime
# number format on col E) so the new row matches the sheet's row template.
$ws.Range("A124:V124").Copy()
$ws.Range("A125:V125").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A125").Value = 124
$ws.Range("B125").Value = "slovakia"
$ws.Range("C125").Value = "2-liga"
$ws.Range("D125").Value = "2023-2024"
$ws.Range("E125").Value = 45241.75
$ws.Range("F125").Value = "Myjava"
$ws.Range("G125").Value = 3
$ws.Range("H125").Value = "Puchov"
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 2.52
$ws.Range("K125").Value = "10/11/2023 04:42"
$ws.Range("L125").Value = 2.31
$ws.Range("M125").Value = "11/11/2023 17:59"
$ws.Range("N125").Value = 3.33
$ws.Range("O125").Value = "10/11/2023 04:42"
$ws.Range("P125").Value = 3.56
$ws.Range("Q125").Value = "11/11/2023 17:59"
$ws.Range("R125").Value = 2.4
$ws.Range("S125").Value = "10/11/2023 04:42"
$ws.Range("T125").Value = 2.8
$ws.Range("U125").Value = "11/11/2023 17:59"
$ws.Range("V125").Value = "https://www.betexplorer.com/football/slovakia/2-liga/myjava-msk-puchov/KlHurxS2/"
